$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.179.44'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '2.892.18'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.66%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0857'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").Value = '3.346.83'
$ws.Range("E15").Value = '  +3.50%  '
$ws.Range("D16").Value = '2.924.43'
$ws.Range("E16").Value = '  +4.45%  '
$ws.Range("E17").Value = '  +6.43%  '
$ws.Range("D18").Value = '52.222.71'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.08%  '
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.29%  '
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("E32").Value = '  +3.80%  '
$ws.Range("E33").Value = '  +8.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '53.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0944'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0461'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.73%  '
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("E38").Value = '  +6.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("E40").Value = '  +3.62%  '
$ws.Range("E41").Value = '  +5.69%  '
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("E43").Value = '  +3.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '121.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.89%  '
$ws.Range("D47").Value = '2.209.19'
$ws.Range("E47").Value = '  +3.54%  '
$ws.Range("E48").Value = '  +6.06%  '
$ws.Range("E49").Value = '  +19.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.950'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.45%  '
